# Updated DPM integration testfixture with hierarchy node labels
#
# Regenerates the UUID "ID" values (column A) on each data sheet of the
# typed-domains-2018-1 fixture, mirroring a re-run of the fixture/hierarchy
# generator that produced fresh node identifiers.

$wb = $excel.ActiveWorkbook

# --- CodeSchemes -----------------------------------------------------
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Range("A2").Value = "2299a20c-3bb2-41d5-b158-b2d9b24511db"

# --- Codes -------------------------------------------------------------
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A2").Value = "87ca86c4-df8c-4660-9be7-047a0add46f0"
$wsCodes.Range("A3").Value = "ea4363c5-5162-4737-bd03-1e7565d3c524"
$wsCodes.Range("A4").Value = "2bf2755f-7501-4626-a8a8-fe7ab788f422"
$wsCodes.Range("A5").Value = "a74070df-d59e-4e19-87b8-d720b96852c8"
$wsCodes.Range("A6").Value = "f1392d9b-1deb-4961-8d47-e486e3205e8c"
$wsCodes.Range("A7").Value = "46efc314-6b14-49fe-aa0e-924137d83737"
$wsCodes.Range("A8").Value = "7a163caa-eb9e-466f-8c27-1c2edd9b995e"
$wsCodes.Range("A9").Value = "127170e7-cf38-4571-9972-709c856d7417"
# Column A re-measures slightly narrower for the new id set.
$wsCodes.Columns.Item(1).ColumnWidth = 31.1428571428571

# --- Extensions ----------------------------------------------------------
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Range("A2").Value = "4eac6323-a1b0-400f-9c63-b4503a058690"

# --- Members_dpmTypedDomain ---------------------------------------------
$wsMembers = $wb.Worksheets.Item("Members_dpmTypedDomain")
$wsMembers.Range("A2").Value = "957db9f9-0959-4e0b-b854-1d6caf5153e0"
$wsMembers.Range("A3").Value = "c55dd8f7-55e2-4b4a-923d-9b8aaf9fc960"
$wsMembers.Range("A4").Value = "31e02f70-cc31-4255-a520-95f31d4c62c1"
$wsMembers.Range("A5").Value = "2706a0b6-23f4-4f67-b42c-9b1b1aa511e3"
$wsMembers.Range("A6").Value = "476253e6-5aa4-4a80-9c06-bff2520dcc29"
$wsMembers.Range("A7").Value = "bfc3f817-411d-4cda-8320-28c801400b6c"
$wsMembers.Range("A8").Value = "a7734bbc-2aa6-4153-8822-006f5e4468c8"
$wsMembers.Range("A9").Value = "279332fe-7d82-4b9a-985b-b7fa5998fe00"
# Column A re-measures slightly narrower for the new id set.
$wsMembers.Columns.Item(1).ColumnWidth = 33.4285714285714
